$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$shortnames = @(
    "shortname",
    "ADIDAS AG",
    "AIRBUS SE",
    "ALLIANZ SE",
    "BASF SE",
    "BAYER AG",
    "BEIERSDORF AG O.N.",
    "BAYERISCHE MOTOREN WERKE AG",
    "BRENNTAG SE NA O.N.",
    "CONTINENTAL AG",
    "COVESTRO AG",
    "Delivery Hero SE Namens-Aktien",
    "DEUTSCHE BANK AG",
    "DEUTSCHE BOERSE NA O.N.",
    "DEUTSCHE POST AG",
    "DEUTSCHE TELEKOM AG",
    "E.ON SE",
    "FRESENIUS MEDICAL CARE AG & CO",
    "FRESENIUS SE&CO KGAA",
    "HEIDELBERGCEMENT AG O.N.",
    "HELLOFRESH SE INH O.N.",
    "HENKEL AG&CO. KGAA",
    "INFINEON TECHNOLOGIES AG",
    "LINDE PLC EO 0,001",
    "MERCEDES-BENZ GROUP",
    "MERCK KGAA",
    "MTU AERO ENGINES NA O.N.",
    "MUENCHENER RUECKVERSICHERUNGS A",
    "PORSCHE AUTOM.HLDG VZO",
    "PUMA SE",
    "Qiagen N.V.",
    "RWE AG INH O.N.",
    "SAP SE",
    "SARTORIUS AG VZO O.N.",
    "SIEMENS AG",
    "SIEMENS ENERGY AG NA O.N.",
    "SIEMENS HEALTH.AG NA O.N.",
    "SYMRISE AG INH. O.N.",
    "VOLKSWAGEN AG",
    "VONOVIA SE NA O.N.",
    "ZALANDO SE"
)

for ($i = 0; $i -lt $shortnames.Length; $i++) {
    $ws.Cells.Item($i + 1, 6).Value = $shortnames[$i]
}

$ws.Range("F1").Font.Bold = $true

$ws.Application.Goto($ws.Range("H46"))
$ws.Application.ActiveWindow.ScrollRow = 12
$ws.Application.ActiveWindow.ScrollColumn = 1
